# Scheduled runner update: refresh cached Universalis market-price derived
# figures (currentAveragePrice*, LevePrice*, LeveProfit*) across the Leve
# crafting-profit sheets. Values below are the latest pulled averages.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2103.25
$ws.Range("I15").Value = 2103.25
$ws.Range("K15").Value = 6309.75
$ws.Range("M15").Value = -6140.75
$ws.Range("H112").Value = 42003.76
$ws.Range("J112").Value = 85091.5
$ws.Range("L112").Value = 255274.5
$ws.Range("N112").Value = -257490.5
$ws.Range("H132").Value = 1994.4546
$ws.Range("I132").Value = 1884.25
$ws.Range("K132").Value = 5652.75
$ws.Range("M132").Value = -3122.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18519630
$ws.Range("I32").Value = 20001082
$ws.Range("K32").Value = 20001082
$ws.Range("M32").Value = -20000795
$ws.Range("H45").Value = 3317.9656
$ws.Range("I45").Value = 2379.9092
$ws.Range("K45").Value = 2379.9092
$ws.Range("M45").Value = -2002.9092
$ws.Range("H55").Value = 50021
$ws.Range("I55").Value = 30048
$ws.Range("J55").Value = 69994
$ws.Range("K55").Value = 30048
$ws.Range("L55").Value = 69994
$ws.Range("M55").Value = -29733
$ws.Range("N55").Value = -70624
$ws.Range("H102").Value = 2138.5833
$ws.Range("I102").Value = 1816.3
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 1816.3
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -194.3
$ws.Range("N102").Value = -6994
$ws.Range("H110").Value = 2263.3684
$ws.Range("I110").Value = 2088.4707
$ws.Range("K110").Value = 2088.4707
$ws.Range("M110").Value = -43.47069999999985
$ws.Range("H122").Value = 6635.76
$ws.Range("I122").Value = 5927.857
$ws.Range("K122").Value = 17783.571
$ws.Range("M122").Value = -15333.571
$ws.Range("H123").Value = 96994.5
$ws.Range("J123").Value = 96994.5
$ws.Range("L123").Value = 96994.5
$ws.Range("N123").Value = -106794.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2559.1538
$ws.Range("I105").Value = 1712.5
$ws.Range("J105").Value = 3284.8572
$ws.Range("K105").Value = 1712.5
$ws.Range("L105").Value = 3284.8572
$ws.Range("M105").Value = 34.5
$ws.Range("N105").Value = -6778.8572
$ws.Range("H107").Value = 1749.5
$ws.Range("I107").Value = 999.5
$ws.Range("J107").Value = 2499.5
$ws.Range("K107").Value = 999.5
$ws.Range("L107").Value = 2499.5
$ws.Range("M107").Value = 920.5
$ws.Range("N107").Value = -6339.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1594155.9
$ws.Range("I6").Value = 1821749.8
$ws.Range("J6").Value = 999
$ws.Range("K6").Value = 1821749.8
$ws.Range("L6").Value = 999
$ws.Range("M6").Value = -1821636.8
$ws.Range("N6").Value = -1225
$ws.Range("H31").Value = 3751.016
$ws.Range("I31").Value = 1430.0605
$ws.Range("J31").Value = 6392.1035
$ws.Range("K31").Value = 1430.0605
$ws.Range("L31").Value = 6392.1035
$ws.Range("M31").Value = -1135.0605
$ws.Range("N31").Value = -6982.1035
$ws.Range("H34").Value = 3751.016
$ws.Range("I34").Value = 1430.0605
$ws.Range("J34").Value = 6392.1035
$ws.Range("K34").Value = 1430.0605
$ws.Range("L34").Value = 6392.1035
$ws.Range("M34").Value = -1228.0605
$ws.Range("N34").Value = -6796.1035
$ws.Range("H81").Value = 43329.332
$ws.Range("J81").Value = 39994.5
$ws.Range("L81").Value = 39994.5
$ws.Range("N81").Value = -41990.5
$ws.Range("H84").Value = 43329.332
$ws.Range("J84").Value = 39994.5
$ws.Range("L84").Value = 119983.5
$ws.Range("N84").Value = -129967.5
$ws.Range("H87").Value = 81454
$ws.Range("J87").Value = 81454
$ws.Range("L87").Value = 81454
$ws.Range("N87").Value = -83826
$ws.Range("H88").Value = 44585.25
$ws.Range("J88").Value = 44585.25
$ws.Range("L88").Value = 44585.25
$ws.Range("N88").Value = -45397.25
$ws.Range("H90").Value = 81454
$ws.Range("J90").Value = 81454
$ws.Range("L90").Value = 244362
$ws.Range("N90").Value = -256218
$ws.Range("H91").Value = 44585.25
$ws.Range("J91").Value = 44585.25
$ws.Range("L91").Value = 44585.25
$ws.Range("N91").Value = -47393.25
$ws.Range("H107").Value = 77596.53999999999
$ws.Range("I107").Value = 83779.586
$ws.Range("K107").Value = 83779.586
$ws.Range("M107").Value = -81859.586
$ws.Range("H114").Value = 56497.5
$ws.Range("J114").Value = 56497.5
$ws.Range("L114").Value = 56497.5
$ws.Range("N114").Value = -65175.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 58
$ws.Range("I8").Value = 58
$ws.Range("K8").Value = 174
$ws.Range("M8").Value = -35
$ws.Range("H14").Value = 9545.666999999999
$ws.Range("I14").Value = 9545.666999999999
$ws.Range("K14").Value = 28637.001
$ws.Range("M14").Value = -28464.001
$ws.Range("H38").Value = 275.7619
$ws.Range("J38").Value = 767.8570999999999
$ws.Range("L38").Value = 2303.5713
$ws.Range("N38").Value = -2997.5713
$ws.Range("H97").Value = 344.75
$ws.Range("I97").Value = 293.33334
$ws.Range("K97").Value = 880.0000200000001
$ws.Range("M97").Value = -384.0000200000001
$ws.Range("H98").Value = 1674
$ws.Range("I98").Value = 2158.6
$ws.Range("K98").Value = 6475.799999999999
$ws.Range("M98").Value = -4977.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4508.5
$ws.Range("I43").Value = 4508.5
$ws.Range("K43").Value = 4508.5
$ws.Range("M43").Value = -4357.5
$ws.Range("H95").Value = 99975.336
$ws.Range("J95").Value = 99975.336
$ws.Range("L95").Value = 99975.336
$ws.Range("N95").Value = -105467.336
$ws.Range("H97").Value = 635.2917
$ws.Range("I97").Value = 504.85715
$ws.Range("J97").Value = 1548.3334
$ws.Range("K97").Value = 504.85715
$ws.Range("L97").Value = 1548.3334
$ws.Range("M97").Value = -8.85714999999999
$ws.Range("N97").Value = -2540.3334
$ws.Range("H107").Value = 1028.8422
$ws.Range("I107").Value = 1232.875
$ws.Range("J107").Value = 880.4545000000001
$ws.Range("K107").Value = 1232.875
$ws.Range("L107").Value = 880.4545000000001
$ws.Range("M107").Value = 687.125
$ws.Range("N107").Value = -4720.4545
$ws.Range("H119").Value = 62011
$ws.Range("J119").Value = 62011
$ws.Range("L119").Value = 62011
$ws.Range("N119").Value = -71687
$ws.Range("H122").Value = 2056
$ws.Range("I122").Value = 2100.6667
$ws.Range("J122").Value = 1966.6666
$ws.Range("K122").Value = 6302.000100000001
$ws.Range("L122").Value = 5899.9998
$ws.Range("M122").Value = -3852.000100000001
$ws.Range("N122").Value = -10799.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3423.2856
$ws.Range("I16").Value = 3160.5
$ws.Range("K16").Value = 3160.5
$ws.Range("M16").Value = -2990.5
$ws.Range("H68").Value = 1668.2858
$ws.Range("I68").Value = 1445.8334
$ws.Range("J68").Value = 3003
$ws.Range("K68").Value = 1445.8334
$ws.Range("L68").Value = 3003
$ws.Range("M68").Value = -696.8334
$ws.Range("N68").Value = -4501
$ws.Range("H71").Value = 1668.2858
$ws.Range("I71").Value = 1445.8334
$ws.Range("J71").Value = 3003
$ws.Range("K71").Value = 7229.166999999999
$ws.Range("L71").Value = 15015
$ws.Range("M71").Value = -3485.166999999999
$ws.Range("N71").Value = -22503
$ws.Range("H122").Value = 11206.518
$ws.Range("I122").Value = 10610.167
$ws.Range("J122").Value = 12182.363
$ws.Range("K122").Value = 31830.501
$ws.Range("L122").Value = 36547.089
$ws.Range("M122").Value = -29380.501
$ws.Range("N122").Value = -41447.089
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8768.154
$ws.Range("I122").Value = 10757.333
$ws.Range("J122").Value = 7063.143
$ws.Range("K122").Value = 32271.999
$ws.Range("L122").Value = 21189.429
$ws.Range("M122").Value = -29821.999
$ws.Range("N122").Value = -26089.429
